$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of header cell -> new value (only the ones that changed)
$renames = @{
    "D1" = "MI_SA_ET15_2015"
    "E1" = "MI_SA_PCDec_ET15_2015"
    "F1" = "MI_SA_PCDD_ET15_2015"
    "G1" = "MI_SA_ETLow_2015"
    "H1" = "MI_SA_PCDec_ETLow_2015"
    "I1" = "MI_SA_PCDD_ETLow_2015"

    "M1" = "MI_SA_ET15_2070"
    "N1" = "MI_SA_PCDec_ET15_2070"
    "O1" = "MI_SA_PCDD_ET15_2070"
    "P1" = "MI_SA_ETLow_2070"
    "Q1" = "MI_SA_PCDec_ETLow_2070"
    "R1" = "MI_SA_PCDD_ETLow_2070"

    "V1" = "MI_SA_ET15_2150"
    "W1" = "MI_SA_PCDec_ET15_2150"
    "X1" = "MI_SA_PCDD_ET15_2150"
    "Y1" = "MI_SA_ETLow_2150"
    "Z1" = "MI_SA_PCDec_ETLow_2150"
    "AA1" = "MI_SA_PCDD_ETLow_2150"
}

foreach ($addr in $renames.Keys) {
    $ws.Range($addr).Value = $renames[$addr]
}
